$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp caption in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 00:33"

# --- Update per-country statistics (columns B..H) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 3410741
$ws.Range("C4").Value = 55095
$ws.Range("D4").Value = 1515674
$ws.Range("E4").Value = 1757302
$ws.Range("G4").Value = 363
$ws.Range("H4").Value = 137765

# Row 19: Alemania
$ws.Range("B19").Value = 199950
$ws.Range("C19").Value = 138
$ws.Range("E19").Value = 6216

# Row 22: Catar
$ws.Range("B22").Value = 150445
$ws.Range("C22").Value = 5083
$ws.Range("D22").Value = 63451
$ws.Range("E22").Value = 81687
$ws.Range("G22").Value = 188
$ws.Range("H22").Value = 5307

# Row 25: Argentina
$ws.Range("B25").Value = 100166
$ws.Range("C25").Value = 2657
$ws.Range("E25").Value = 55627
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 1845

# Row 34: Kazajistan
$ws.Range("E34").Value = 24064
$ws.Range("G34").Value = 21
$ws.Range("H34").Value = 375

# Row 54: Guatemala
$ws.Range("B54").Value = 29355
$ws.Range("C54").Value = 757
$ws.Range("D54").Value = 4214
$ws.Range("E54").Value = 23922
$ws.Range("G54").Value = 47
$ws.Range("H54").Value = 1219

# Row 57: Ghana
$ws.Range("B57").Value = 24518
$ws.Range("C57").Value = 270
$ws.Range("D57").Value = 20187
$ws.Range("E57").Value = 4192
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 139

# Row 86: Etiopia
$ws.Range("B86").Value = 7252
$ws.Range("C86").Value = 77
$ws.Range("D86").Value = 3319
$ws.Range("E86").Value = 3665
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 268

# --- Paraguay overtakes Nicaragua in total cases, so the two countries swap places ---
# Row 106 becomes Paraguay with its freshly updated figures
$ws.Range("A106").Value = "Paraguay"
$ws.Range("B106").Value = 2948
$ws.Range("C106").Value = 128
$ws.Range("D106").Value = 1275
$ws.Range("E106").Value = 1651
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 22

# Row 107 becomes Nicaragua, carrying the figures that used to sit in row 106
$ws.Range("A107").Value = "Nicaragua"
$ws.Range("B107").Value = 2846
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 1993
$ws.Range("E107").Value = 762
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 91
